# Append the two new journal paragraphs (week 2 continuation) to the end of the document body,
# right before the final paragraph mark, preserving all existing content untouched.
$d = $word.ActiveDocument
$lastPara = $d.Paragraphs.Last
$bodyEnd = $lastPara.Range.End
$insertionPoint = $d.Range($bodyEnd - 1, $bodyEnd - 1)

$newParagraphsXml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t xml:space="preserve">This week I </w:t></w:r><w:r><w:t>used</w:t></w:r><w:r><w:t xml:space="preserve"> the p5 reference page and went to the typography section</w:t></w:r><w:r><w:t xml:space="preserve">, while I </w:t></w:r><w:r><w:t>scrolled,</w:t></w:r><w:r><w:t xml:space="preserve"> I saw the camera section and after thinking about 3d camera I decided to check the camera section out as well. The </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>textalign</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">() tool reminded me of html but I also recall using it for one of the worksheets. </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>t</w:t></w:r><w:r><w:t>extleading</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>() s</w:t></w:r><w:r><w:t xml:space="preserve">ets the number space in between each line of text. </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>textAscent</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">() and </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>textDescent</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">() shows the distance of the tallest character above and below the baseline. </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>loadFont</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">() will be a useful tool as it allows for the import of fonts. p5 showed the syntax to be </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>loadFont</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>(path,[</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>callback</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>,[</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>onError</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>])</w:t></w:r><w:r><w:t xml:space="preserve">. This path is the file or URL to load for the font. The other two were optional. </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>onError</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> seemed to be the more useful tool as it works like a backup plan and will execute another function if an error were to occur. </w:t></w:r><w:r><w:t xml:space="preserve">The text() reference page had an interesting example of rotation being used. It used </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>rotateZ</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">() which gave the text a 3d effect as it rotated on the z axis. </w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t xml:space="preserve">The camera section was more complicated to get my head around and I am unsure if I have. </w:t></w:r><w:r><w:t xml:space="preserve">The camera() tool had one example of a square moving back and forward in the screen – looked to be on a z axis rather than scaling the square up and down. The other example involved sliders. The sliders controlled 6 parameters of the cube. Depending on how you moved the sliders the camera would follow and show you the cube from different angles. The last thing I looked at was perspective(). perspective() involves a 3d sketch which I assume has a static camera in the middle you rotate </w:t></w:r><w:r><w:lastRenderedPageBreak/><w:t>around. Using this</w:t></w:r><w:r><w:t>,</w:t></w:r><w:r><w:t xml:space="preserve"> any image </w:t></w:r><w:r><w:t>closest</w:t></w:r><w:r><w:t xml:space="preserve"> to the perspective</w:t></w:r><w:r><w:t xml:space="preserve"> you are looking at appears larger while any further appear to be in the background. It adds a range of depth to the image.</w:t></w:r></w:p>
'@

$insertionPoint.InsertXML($newParagraphsXml)
